$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-29 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-30 Monday", 2)

$d.Content.Find.Execute("96×66=6336", $true, $false, $false, $false, $false, $true, 1, $false, "81×85=6885", 2)
$d.Content.Find.Execute("20×34=680", $true, $false, $false, $false, $false, $true, 1, $false, "85×95=8075", 2)
$d.Content.Find.Execute("71×75=5325", $true, $false, $false, $false, $false, $true, 1, $false, "90×89=8010", 2)
$d.Content.Find.Execute("44×63=2772", $true, $false, $false, $false, $false, $true, 1, $false, "35×99=3465", 2)
$d.Content.Find.Execute("32×80=2560", $true, $false, $false, $false, $false, $true, 1, $false, "40×79=3160", 2)
$d.Content.Find.Execute("14×76=1064", $true, $false, $false, $false, $false, $true, 1, $false, "96×88=8448", 2)
$d.Content.Find.Execute("20×23=460", $true, $false, $false, $false, $false, $true, 1, $false, "51×73=3723", 2)
$d.Content.Find.Execute("13×60=780", $true, $false, $false, $false, $false, $true, 1, $false, "77×65=5005", 2)
$d.Content.Find.Execute("94×22=2068", $true, $false, $false, $false, $false, $true, 1, $false, "64×45=2880", 2)
$d.Content.Find.Execute("71×45=3195", $true, $false, $false, $false, $false, $true, 1, $false, "67×53=3551", 2)
$d.Content.Find.Execute("56×57=3192", $true, $false, $false, $false, $false, $true, 1, $false, "12×83=996", 2)
$d.Content.Find.Execute("11×87=957", $true, $false, $false, $false, $false, $true, 1, $false, "36×80=2880", 2)
$d.Content.Find.Execute("91×61=5551", $true, $false, $false, $false, $false, $true, 1, $false, "42×97=4074", 2)
$d.Content.Find.Execute("98×54=5292", $true, $false, $false, $false, $false, $true, 1, $false, "32×42=1344", 2)
$d.Content.Find.Execute("71×41=2911", $true, $false, $false, $false, $false, $true, 1, $false, "58×29=1682", 2)
$d.Content.Find.Execute("22×98=2156", $true, $false, $false, $false, $false, $true, 1, $false, "61×15=915", 2)
$d.Content.Find.Execute("86×71=6106", $true, $false, $false, $false, $false, $true, 1, $false, "35×74=2590", 2)
$d.Content.Find.Execute("90×76=6840", $true, $false, $false, $false, $false, $true, 1, $false, "92×57=5244", 2)
$d.Content.Find.Execute("52×11=572", $true, $false, $false, $false, $false, $true, 1, $false, "76×20=1520", 2)
$d.Content.Find.Execute("25×25=625", $true, $false, $false, $false, $false, $true, 1, $false, "34×81=2754", 2)
$d.Content.Find.Execute("44×40=1760", $true, $false, $false, $false, $false, $true, 1, $false, "87×64=5568", 2)
$d.Content.Find.Execute("68×23=1564", $true, $false, $false, $false, $false, $true, 1, $false, "51×58=2958", 2)
$d.Content.Find.Execute("83×31=2573", $true, $false, $false, $false, $false, $true, 1, $false, "48×96=4608", 2)
$d.Content.Find.Execute("82×25=2050", $true, $false, $false, $false, $false, $true, 1, $false, "35×23=805", 2)
$d.Content.Find.Execute("40×90=3600", $true, $false, $false, $false, $false, $true, 1, $false, "59×21=1239", 2)
